$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the logo placeholder cells from row 20 up into row 14 (H14, K14),
# matching the target layout where Controller/H14/K14 share row 14.
$ws.Range("H20:K20").Cut($ws.Range("H14"))

# Tidy up the now-empty source cells so they fully reset to the default,
# unstyled state (no leftover number formatting from the old logo cells).
$ws.Range("H20").NumberFormat = "General"
$ws.Range("K20").NumberFormat = "General"
$ws.Range("H20").ClearContents()
$ws.Range("K20").ClearContents()

# Reflect the resulting selection, as captured in the saved workbook.
$ws.Range("H14:K14").Select()
